$d = $word.ActiveDocument

$d.Content.Find.Execute("96×96=", $true, $false, $false, $false, $false, $true, 1, $false, "71×32=", 2)
$d.Content.Find.Execute("32×71=", $true, $false, $false, $false, $false, $true, 1, $false, "97×48=", 2)
$d.Content.Find.Execute("94×55=", $true, $false, $false, $false, $false, $true, 1, $false, "49×50=", 2)
$d.Content.Find.Execute("18×94=", $true, $false, $false, $false, $false, $true, 1, $false, "69×23=", 2)
$d.Content.Find.Execute("57×49=", $true, $false, $false, $false, $false, $true, 1, $false, "42×52=", 2)
$d.Content.Find.Execute("96×80=", $true, $false, $false, $false, $false, $true, 1, $false, "63×72=", 2)
$d.Content.Find.Execute("39×30=", $true, $false, $false, $false, $false, $true, 1, $false, "81×15=", 2)
$d.Content.Find.Execute("84×99=", $true, $false, $false, $false, $false, $true, 1, $false, "42×51=", 2)
$d.Content.Find.Execute("86×66=", $true, $false, $false, $false, $false, $true, 1, $false, "80×56=", 2)
$d.Content.Find.Execute("17×83=", $true, $false, $false, $false, $false, $true, 1, $false, "29×50=", 2)
$d.Content.Find.Execute("69×86=", $true, $false, $false, $false, $false, $true, 1, $false, "37×23=", 2)
$d.Content.Find.Execute("53×48=", $true, $false, $false, $false, $false, $true, 1, $false, "38×44=", 2)
$d.Content.Find.Execute("20×78=", $true, $false, $false, $false, $false, $true, 1, $false, "20×69=", 2)
$d.Content.Find.Execute("48×64=", $true, $false, $false, $false, $false, $true, 1, $false, "83×34=", 2)
$d.Content.Find.Execute("32×47=", $true, $false, $false, $false, $false, $true, 1, $false, "89×94=", 2)
$d.Content.Find.Execute("87×90=", $true, $false, $false, $false, $false, $true, 1, $false, "74×89=", 2)
$d.Content.Find.Execute("36×93=", $true, $false, $false, $false, $false, $true, 1, $false, "27×92=", 2)
$d.Content.Find.Execute("77×94=", $true, $false, $false, $false, $false, $true, 1, $false, "18×18=", 2)
$d.Content.Find.Execute("15×30=", $true, $false, $false, $false, $false, $true, 1, $false, "77×52=", 2)
$d.Content.Find.Execute("13×15=", $true, $false, $false, $false, $false, $true, 1, $false, "22×48=", 2)
$d.Content.Find.Execute("56×64=", $true, $false, $false, $false, $false, $true, 1, $false, "99×33=", 2)
$d.Content.Find.Execute("94×16=", $true, $false, $false, $false, $false, $true, 1, $false, "51×96=", 2)
$d.Content.Find.Execute("54×34=", $true, $false, $false, $false, $false, $true, 1, $false, "45×62=", 2)
$d.Content.Find.Execute("76×58=", $true, $false, $false, $false, $false, $true, 1, $false, "58×73=", 2)
$d.Content.Find.Execute("40×79=", $true, $false, $false, $false, $false, $true, 1, $false, "56×73=", 2)
